# Regenerate orders with updated distance/size codes.
#   Distance: D80 -> D86, D51 -> D55, D64 -> D69
#   Size:     S30 -> S31   (S20 / S25 unchanged)
# These codes appear embedded inside several text columns (Condition,
# Filename_Left, Filename_Right, Distance, Size) throughout the sheet,
# so we sweep the whole used range and rewrite any text cell that
# contains one of the old tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

$lastRow = $startRow + $rowCount - 1
$lastCol = $startCol + $colCount - 1

for ($r = $startRow; $r -le $lastRow; $r++) {
    for ($c = $startCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -is [string]) {
            if (($v.Contains("D80")) -or ($v.Contains("D51")) -or ($v.Contains("D64")) -or ($v.Contains("S30"))) {
                $nv = $v.Replace("D80", "D86").Replace("D51", "D55").Replace("D64", "D69").Replace("S30", "S31")
                $cell.Value = $nv
            }
        }
    }
}
